$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.951.00"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.674.99"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.03"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.22"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "1.910.96"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.680.36"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.64"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "26.963.39"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.07"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.91"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.16"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.61"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.05"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "1.479.10"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.67"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.584"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.896"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  +7.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.84"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.93"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "1.815.84"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.779"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.50"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.55%  "
